$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "price" cells are plain text in the source (t="inlineStr"),
# but several values look numeric (e.g. "1.001", "316.41"). Force the
# cell to text format before assigning so Excel keeps it as a string
# instead of silently coercing it to a number, then restore the default
# "Normal" style so no stray number-format style sticks to the cell.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "25.931.37"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +8.24%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.764.91"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +6.13%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "316.41"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.55%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9973"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3829"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.98%  "

$ws.Range("E8").Value = "  +5.03%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "50.65"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +6.24%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.233"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +5.32%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07693"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +6.07%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.9977"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.08%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "21.68"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +5.13%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.497"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +7.61%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.101"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +5.16%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.763.61"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +5.98%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001162"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +6.06%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.9977"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06801"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.31%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "87.16"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +6.91%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.75"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +7.97%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.524"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +6.60%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "12.80"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +6.67%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "25.844.48"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +7.90%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.431"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.26%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.952"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +9.51%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.72"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +5.92%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "156.18"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.69%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.963.25"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +6.50%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "134.27"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +5.54%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.228"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +25.33%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "7.242"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +15.03%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.255"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.90%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "14.25"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +14.91%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.812"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +4.57%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.08773"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +4.50%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.725"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +7.83%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.06787"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +6.79%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.02501"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +7.95%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "9.362"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.15%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.2256"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +8.61%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.295"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.65%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.6573"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +7.64%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "14.41"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +9.13%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.9971"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.31%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.6354"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +6.75%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.918"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.52%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.172"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +8.50%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "133.01"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.48%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.07497"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +5.88%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "81.12"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +6.69%  "

